# Refresh currentAveragePrice / LevePrice / LeveProfit figures across the
# Lich_Profits job sheets (scheduled runner data sync).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 477281.8
$ws.Range("J17").Value = 501111.1
$ws.Range("L17").Value = 1503333.3
$ws.Range("N17").Value = -1503669.3
$ws.Range("H43").Value = 1966
$ws.Range("I43").Value = 1900
$ws.Range("K43").Value = 1900
$ws.Range("M43").Value = -1831
$ws.Range("H113").Value = 5367.1787
$ws.Range("I113").Value = 9158
$ws.Range("J113").Value = 3571.5264
$ws.Range("K113").Value = 9158
$ws.Range("L113").Value = 3571.5264
$ws.Range("M113").Value = -5904
$ws.Range("N113").Value = -10079.5264
$ws.Range("H116").Value = 5449.2
$ws.Range("I116").Value = 4309
$ws.Range("J116").Value = 5734.25
$ws.Range("K116").Value = 4309
$ws.Range("L116").Value = 5734.25
$ws.Range("M116").Value = -867
$ws.Range("N116").Value = -12618.25
$ws.Range("H132").Value = 4716.3
$ws.Range("I132").Value = 1984.8823
$ws.Range("K132").Value = 5954.6469
$ws.Range("M132").Value = -3424.6469
$ws.Range("H136").Value = 267497
$ws.Range("J136").Value = 267497
$ws.Range("L136").Value = 267497
$ws.Range("N136").Value = -277697

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 9157.842000000001
$ws.Range("I6").Value = 9111.056
$ws.Range("K6").Value = 9111.056
$ws.Range("M6").Value = -8938.056
$ws.Range("H32").Value = 5786.8228
$ws.Range("I32").Value = 3922.2957
$ws.Range("J32").Value = 22334.5
$ws.Range("K32").Value = 3922.2957
$ws.Range("L32").Value = 22334.5
$ws.Range("M32").Value = -3635.2957
$ws.Range("N32").Value = -22908.5
$ws.Range("H39").Value = 1608
$ws.Range("I39").Value = 1608
$ws.Range("K39").Value = 1608
$ws.Range("M39").Value = -1088
$ws.Range("H61").Value = 10629.846
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 15711
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 15711
$ws.Range("M61").Value = -2288
$ws.Range("N61").Value = -16135
$ws.Range("H136").Value = 10629.846
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 15711
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 47133
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -52233
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H94").Value = 812.58826
$ws.Range("I94").Value = 858.3214
$ws.Range("K94").Value = 858.3214
$ws.Range("M94").Value = -407.3214
$ws.Range("H99").Value = 4617.278
$ws.Range("I99").Value = 4703.125
$ws.Range("J99").Value = 3930.5
$ws.Range("K99").Value = 4703.125
$ws.Range("L99").Value = 3930.5
$ws.Range("M99").Value = -3205.125
$ws.Range("N99").Value = -6926.5
$ws.Range("H126").Value = 45000
$ws.Range("J126").Value = 45000
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -54880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1114187.8
$ws.Range("I31").Value = 1114187.8
$ws.Range("K31").Value = 1114187.8
$ws.Range("M31").Value = -1113892.8
$ws.Range("H34").Value = 1114187.8
$ws.Range("I34").Value = 1114187.8
$ws.Range("K34").Value = 1114187.8
$ws.Range("M34").Value = -1113985.8
$ws.Range("H68").Value = 37761.5
$ws.Range("J68").Value = 39713.8
$ws.Range("L68").Value = 39713.8
$ws.Range("N68").Value = -41211.8
$ws.Range("H71").Value = 37761.5
$ws.Range("J71").Value = 39713.8
$ws.Range("L71").Value = 119141.4
$ws.Range("N71").Value = -126629.4
$ws.Range("H74").Value = 99313.5
$ws.Range("J74").Value = 99313.5
$ws.Range("L74").Value = 99313.5
$ws.Range("N74").Value = -101061.5
$ws.Range("H77").Value = 99313.5
$ws.Range("J77").Value = 99313.5
$ws.Range("L77").Value = 297940.5
$ws.Range("N77").Value = -306676.5
$ws.Range("H94").Value = 1447.3334
$ws.Range("I94").Value = 1106.5
$ws.Range("J94").Value = 1571.2727
$ws.Range("K94").Value = 1106.5
$ws.Range("L94").Value = 1571.2727
$ws.Range("M94").Value = -655.5
$ws.Range("N94").Value = -2473.2727
$ws.Range("H99").Value = 1278413.6
$ws.Range("I99").Value = 3356670.8
$ws.Range("J99").Value = 31459.4
$ws.Range("K99").Value = 3356670.8
$ws.Range("L99").Value = 31459.4
$ws.Range("M99").Value = -3355172.8
$ws.Range("N99").Value = -34455.4
$ws.Range("H109").Value = 35277.5
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27080
$ws.Range("H126").Value = 1278413.6
$ws.Range("I126").Value = 3356670.8
$ws.Range("J126").Value = 31459.4
$ws.Range("K126").Value = 10070012.4
$ws.Range("L126").Value = 94378.20000000001
$ws.Range("M126").Value = -10067542.4
$ws.Range("N126").Value = -99318.20000000001
$ws.Range("H134").Value = 5604.9585
$ws.Range("I134").Value = 5919.952
$ws.Range("K134").Value = 17759.856
$ws.Range("M134").Value = -15224.856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6665
$ws.Range("J80").Value = 6968.875
$ws.Range("L80").Value = 20906.625
$ws.Range("N80").Value = -22778.625
$ws.Range("H83").Value = 6665
$ws.Range("J83").Value = 6968.875
$ws.Range("L83").Value = 62719.875
$ws.Range("N83").Value = -72079.875
$ws.Range("H129").Value = 3289.3
$ws.Range("I129").Value = 718
$ws.Range("K129").Value = 2154
$ws.Range("M129").Value = 2846
$ws.Range("H141").Value = 3426.0908
$ws.Range("I141").Value = 3187.4443
$ws.Range("K141").Value = 9562.332900000001
$ws.Range("M141").Value = -4382.332900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 13635.143
$ws.Range("I126").Value = 13635.143
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 40905.429
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -38435.429
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3314.0667
$ws.Range("I22").Value = 3005.875
$ws.Range("K22").Value = 3005.875
$ws.Range("M22").Value = -2710.875
$ws.Range("H27").Value = 3314.0667
$ws.Range("I27").Value = 3005.875
$ws.Range("K27").Value = 3005.875
$ws.Range("M27").Value = -2898.875
$ws.Range("H41").Value = 31480.8
$ws.Range("I41").Value = 30007.75
$ws.Range("J41").Value = 37373
$ws.Range("K41").Value = 30007.75
$ws.Range("L41").Value = 37373
$ws.Range("M41").Value = -29569.75
$ws.Range("N41").Value = -38249
$ws.Range("H132").Value = 5461.778
$ws.Range("I132").Value = 5767.1353
$ws.Range("J132").Value = 4797.1763
$ws.Range("K132").Value = 17301.4059
$ws.Range("L132").Value = 14391.5289
$ws.Range("M132").Value = -14771.4059
$ws.Range("N132").Value = -19451.5289
$ws.Range("H136").Value = 3185.0513
$ws.Range("I136").Value = 1885.3636
$ws.Range("J136").Value = 10333.333
$ws.Range("K136").Value = 5656.0908
$ws.Range("L136").Value = 30999.999
$ws.Range("M136").Value = -3106.0908
$ws.Range("N136").Value = -36099.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1243.8636
$ws.Range("I113").Value = 859.7059
$ws.Range("K113").Value = 2579.1177
$ws.Range("M113").Value = -409.1177000000002
$ws.Range("H126").Value = 2721
$ws.Range("I126").Value = 2723.875
$ws.Range("K126").Value = 8171.625
$ws.Range("M126").Value = -5701.625
$ws.Range("H132").Value = 1577.0682
$ws.Range("I132").Value = 1170.1613
$ws.Range("J132").Value = 2547.3845
$ws.Range("K132").Value = 3510.4839
$ws.Range("L132").Value = 7642.1535
$ws.Range("M132").Value = -980.4839000000002
$ws.Range("N132").Value = -12702.1535
$ws.Range("H136").Value = 327353.66
$ws.Range("I136").Value = 375517.25
$ws.Range("K136").Value = 1126551.75
$ws.Range("M136").Value = -1124001.75

